$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<hotel>"
$ws.Range("C2").Value = 45

# Row 3
$ws.Range("B3").Value = "<tan>"
$ws.Range("C3").Value = 50

# Row 4
$ws.Range("B4").Value = "<but>"
$ws.Range("C4").Value = 44

# Row 5
$ws.Range("B5").Value = "<is>"
$ws.Range("C5").Value = 44

# Row 6
$ws.Range("B6").Value = "<was>"
$ws.Range("C6").Value = 45

# Row 7
$ws.Range("C7").Value = 41

# Row 8
$ws.Range("B8").Value = "<made>"
$ws.Range("C8").Value = 43

# Row 9
$ws.Range("B9").Value = "<use>"
$ws.Range("C9").Value = 44

# Row 10
$ws.Range("B10").Value = "<the>"
$ws.Range("C10").Value = 43

# Row 12
$ws.Range("B12").Value = "<bram>"
$ws.Range("C12").Value = 50

# Row 13
$ws.Range("B13").Value = "<four>"
$ws.Range("C13").Value = 37

# Row 14
$ws.Range("C14").Value = 42

# Row 15
$ws.Range("B15").Value = "<a>"
$ws.Range("C15").Value = 51

# Row 16
$ws.Range("B16").Value = "<from>"
$ws.Range("C16").Value = 47

# Row 17
$ws.Range("C17").Value = 50

# Row 18
$ws.Range("C18").Value = 42
